$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2023-10-06 -> 2023-10-07, i.e. 45205 -> 45206) for every data row (2..469).
$ws.Range("C2:C469").Value = 45206
